# edit.ps1 - apply the "Fix typo, pdf update" commit to the deck.
#
# 1) Handout-master & notes-master footer date field: 2020-06-30 -> 2020-07-01
# 2) Slide 12 body text typo fix:
#      "static method 를 사용할 것"  ->  "function을 사용할 것"

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Refresh the cached "datetimeFigureOut" footer field text on the handout
#    master and the notes master (both show the footer date as 2020-06-30).
# ---------------------------------------------------------------------------
$hm = $p.HandoutMaster
$hm.HeadersFooters.DateAndTime.Text = "2020-07-01"

$nm = $p.NotesMaster
$nm.HeadersFooters.DateAndTime.Text = "2020-07-01"

# ---------------------------------------------------------------------------
# 2) Fix the "static method 를 사용할 것" -> "function을 사용할 것" typo on the
#    "Exercise" slide. Locate the shape by its current text instead of a
#    hard-coded index so the script keeps working even if shapes shift.
# ---------------------------------------------------------------------------
$needle = "static method "
$targetRange = $null

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $sl = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $sl.Shapes.Count; $shi++) {
        $sh = $sl.Shapes.Item($shi)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.HasText) {
                $tf = $sh.TextFrame
                $tr = $tf.TextRange
                for ($pi = 1; $pi -le $tr.Paragraphs().Count; $pi++) {
                    $para = $tr.Paragraphs($pi, 1)
                    if ($para.Text.Contains($needle)) {
                        $targetRange = $para
                    }
                }
            }
        }
    }
}

# Replace "static method " (keeps its own run/formatting) with "function".
$paraText = $targetRange.Text
$idx = $paraText.IndexOf($needle)
$run = $targetRange.Characters($idx + 1, $needle.Length)
$run.Text = "function"

# Replace the following "를 " with "을 " as its own run, so the trailing
# "사용할 것" keeps its original, untouched run.
$paraText2 = $targetRange.Text
$needle2 = "를 "
$idx2 = $paraText2.IndexOf($needle2)
$run2 = $targetRange.Characters($idx2 + 1, $needle2.Length)
$run2.Text = "을 "
